# Add a new "Vimar" producer row (Usd, 8.15) below the existing data
# and move the active selection to E9, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Vimar"
$ws.Range("B10").Value = "Usd"
$ws.Range("C10").Value = 8.15

$ws.Range("E9").Select()
